$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- USART calculation header ---
$ws.Range("C1").Value = "USART calculation"
$ws.Range("C1").Font.Bold = $true
$ws.Range("C1").Font.Size = 14
$ws.Rows.Item(1).RowHeight = 18.75

# --- I2C calculation section ---
$ws.Range("C9").Value = "I2C calculation"
$ws.Range("C9").Font.Bold = $true
$ws.Range("C9").Font.Size = 14
$ws.Rows.Item(9).RowHeight = 18.75

$ws.Range("D11").Value = "F_CPU"
$ws.Range("E11").Value = 9216000

$ws.Range("D12").Value = "Prescaler TWPS"
$ws.Range("E12").Value = 0

$ws.Range("D13").Value = "TWI Bit rate register TWBR"
$ws.Range("E13").Value = 39

$ws.Range("D14").Value = "SCL frequency"
$ws.Range("E14").Formula = "= E11 / (16 + (2*E13*(POWER(4,E12))))"

$ws.Range("D16").Value = "set SCL frequency"
$ws.Range("E16").Value = 100000

$ws.Range("D17").Value = "set prescaler"
$ws.Range("E17").Value = 0

$ws.Range("D18").Value = "result TWI bit rate register value"
$ws.Range("E18").Formula = "=FLOOR(((E11/E16)-15)/2,1)"

$ws.Range("D19").Value = "actual SCR frequency"
$ws.Range("E19").Formula = "= E11 / (16 + (2*E18*(POWER(4,E17))))"

# Update selection to match the target state
$ws.Range("E18").Select()
